$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Inner_Outer"
$ws.Range("F1").Value = "marpizza_price"
$ws.Range("G1").Value = "beverage_price"
$ws.Range("H1").Value = "distance2ceu"

$ws.Range("F2").Value = 2300

$ws.Range("H1").Select()
